# Update Quantities worksheet:
#  - Rows 2-30: shift the date in column A forward (values in B-J are unchanged)
#  - Rows 31-40: new rows appended at the bottom holding the "older" dates/values
#    that rolled off the top of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New dates for existing rows 2..30 (column A) ----
$datesExisting = @(45569,45570,45571,45572,45573,45574,45575,45576,45577,45578,45579,45580,45581,45582,45583,45584,45585,45586,45587,45588,45589,45590,45591,45592,45593,45594,45595,45596,45568)

for ($i = 0; $i -lt $datesExisting.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $datesExisting[$i]
}

# Rows 2 and 3 previously held the oldest dates (45558/45559) and therefore
# used the "old" C/G/J value set; now that they hold newer dates they must
# switch to the "current" C/G/J value set used by the rest of the table.
$ws.Range("C2").Value = 0.00170247
$ws.Range("G2").Value = 465.80531254
$ws.Range("J2").Value = 485.38834923

$ws.Range("C3").Value = 0.00170247
$ws.Range("G3").Value = 465.80531254
$ws.Range("J3").Value = 485.38834923

# ---- New rows 31..40 ----
$newDates  = @(45567,45566,45564,45565,45563,45558,45559,45560,45561,45562)
$newC      = @(0.00170247,0.00170247,0.00170247,0.00170247,0.00170247,0.00004012,0.00004012,0.00170247,0.00170247,0.00170247)
$newG      = @(465.80531254,465.80531254,465.80531254,465.80531254,465.80531254,280.99031254,280.99031254,465.80531254,465.80531254,465.80531254)
$newJ      = @(485.38834923,485.38834923,485.38834923,485.38834923,485.38834923,1941.48834923,1941.48834923,485.38834923,485.38834923,485.38834923)

# Common (constant across every data row) values for columns B, D, E, F, H, I
$constB = 116.4121952
$constD = 0.008850780000000001
$constE = 0.06933635
$constF = 12792.90181321
$constH = 0.24
$constI = 1.7904431

# Give the new A-column cells the same formatting (bold, border, centered,
# date number format) as the existing date column by copying the style
# from A2 down onto A31:A40 before writing the values.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A31:A40").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $i + 31
    $ws.Cells.Item($r, 1).Value  = $newDates[$i]
    $ws.Cells.Item($r, 2).Value  = $constB
    $ws.Cells.Item($r, 3).Value  = $newC[$i]
    $ws.Cells.Item($r, 4).Value  = $constD
    $ws.Cells.Item($r, 5).Value  = $constE
    $ws.Cells.Item($r, 6).Value  = $constF
    $ws.Cells.Item($r, 7).Value  = $newG[$i]
    $ws.Cells.Item($r, 8).Value  = $constH
    $ws.Cells.Item($r, 9).Value  = $constI
    $ws.Cells.Item($r, 10).Value = $newJ[$i]
}
